$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the full target dataset (header row + 121 component rows) as a 2D array
$data = New-Object 'object[,]' 122,5
$data[0,0] = "Designator"
$data[0,1] = "Mid X"
$data[0,2] = "Mid Y"
$data[0,3] = "Rotation"
$data[0,4] = "Layer"
$data[1,0] = "C1"
$data[1,1] = 108.7565
$data[1,2] = -130.84169299999999
$data[1,3] = 0
$data[1,4] = "top"
$data[2,0] = "C2"
$data[2,1] = 113.8965
$data[2,2] = -123.24369299999999
$data[2,3] = 0
$data[2,4] = "top"
$data[3,0] = "C3"
$data[3,1] = 113.84650000000001
$data[3,2] = -130.14369300000001
$data[3,3] = 180
$data[3,4] = "top"
$data[4,0] = "C4"
$data[4,1] = 108.7565
$data[4,2] = -128.56839299999999
$data[4,3] = 180
$data[4,4] = "top"
$data[5,0] = "C5"
$data[5,1] = 106.2165
$data[5,2] = -121.18969300000001
$data[5,3] = -90
$data[5,4] = "top"
$data[6,0] = "C6"
$data[6,1] = 108.7565
$data[6,2] = -126.269693
$data[6,3] = 180
$data[6,4] = "top"
$data[7,0] = "C7"
$data[7,1] = 103.93049999999999
$data[7,2] = -121.18969300000001
$data[7,3] = 90
$data[7,4] = "top"
$data[8,0] = "C8"
$data[8,1] = 108.7565
$data[8,2] = -133.12769299999999
$data[8,3] = 180
$data[8,4] = "top"
$data[9,0] = "C9"
$data[9,1] = 108.7565
$data[9,2] = -123.983693
$data[9,3] = 0
$data[9,4] = "top"
$data[10,0] = "C10"
$data[10,1] = 86.846500000000006
$data[10,2] = -119.343693
$data[10,3] = -90
$data[10,4] = "top"
$data[11,0] = "C11"
$data[11,1] = 108.7565
$data[11,2] = -135.41369299999999
$data[11,3] = 180
$data[11,4] = "top"
$data[12,0] = "C12"
$data[12,1] = 96.8185
$data[12,2] = -122.205693
$data[12,3] = 90
$data[12,4] = "top"
$data[13,0] = "C13"
$data[13,1] = 77.396500000000003
$data[13,2] = -125.043693
$data[13,3] = 180
$data[13,4] = "top"
$data[14,0] = "C14"
$data[14,1] = 77.296499999999995
$data[14,2] = -119.24369299999999
$data[14,3] = 0
$data[14,4] = "top"
$data[15,0] = "C15"
$data[15,1] = 147.79650000000001
$data[15,2] = -138.74369300000001
$data[15,3] = -90
$data[15,4] = "top"
$data[16,0] = "C16"
$data[16,1] = 65.6965
$data[16,2] = -140.79369299999999
$data[16,3] = 90
$data[16,4] = "top"
$data[17,0] = "C17"
$data[17,1] = 150.59649999999999
$data[17,2] = -138.54369299999999
$data[17,3] = 90
$data[17,4] = "top"
$data[18,0] = "C18"
$data[18,1] = 148.29650000000001
$data[18,2] = -143.24369300000001
$data[18,3] = 90
$data[18,4] = "top"
$data[19,0] = "C19"
$data[19,1] = 166.69
$data[19,2] = -131.93000000000001
$data[19,3] = 90
$data[19,4] = "top"
$data[20,0] = "C20"
$data[20,1] = 91.484499999999997
$data[20,2] = -126.396693
$data[20,3] = 0
$data[20,4] = "top"
$data[21,0] = "C21"
$data[21,1] = 123.3965
$data[21,2] = -105.343693
$data[21,3] = -90
$data[21,4] = "top"
$data[22,0] = "C22"
$data[22,1] = 69.6965
$data[22,2] = -140.093693
$data[22,3] = -90
$data[22,4] = "top"
$data[23,0] = "C23"
$data[23,1] = 69.6965
$data[23,2] = -143.89369300000001
$data[23,3] = 90
$data[23,4] = "top"
$data[24,0] = "C24"
$data[24,1] = 101.64449999999999
$data[24,2] = -121.18969300000001
$data[24,3] = 90
$data[24,4] = "top"
$data[25,0] = "C25"
$data[25,1] = 91.484499999999997
$data[25,2] = -134.016693
$data[25,3] = 0
$data[25,4] = "top"
$data[26,0] = "C26"
$data[26,1] = 91.738500000000002
$data[26,2] = -137.699693
$data[26,3] = 90
$data[26,4] = "top"
$data[27,0] = "C27"
$data[27,1] = 82.594499999999996
$data[27,2] = -127.793693
$data[27,3] = 0
$data[27,4] = "top"
$data[28,0] = "C28"
$data[28,1] = 102.02549999999999
$data[28,2] = -137.699693
$data[28,3] = 90
$data[28,4] = "top"
$data[29,0] = "C29"
$data[29,1] = 91.484499999999997
$data[29,2] = -131.60369299999999
$data[29,3] = 0
$data[29,4] = "top"
$data[30,0] = "C30"
$data[30,1] = 91.484499999999997
$data[30,2] = -128.80969300000001
$data[30,3] = 180
$data[30,4] = "top"
$data[31,0] = "C31"
$data[31,1] = 104.1845
$data[31,2] = -137.699693
$data[31,3] = 90
$data[31,4] = "top"
$data[32,0] = "C32"
$data[32,1] = 99.358500000000006
$data[32,2] = -142.652693
$data[32,3] = 90
$data[32,4] = "top"
$data[33,0] = "C33"
$data[33,1] = 82.594499999999996
$data[33,2] = -134.14369300000001
$data[33,3] = 180
$data[33,4] = "top"
$data[34,0] = "C34"
$data[34,1] = 194.95500000000001
$data[34,2] = -129.27500000000001
$data[34,3] = 90
$data[34,4] = "top"
$data[35,0] = "C35"
$data[35,1] = 171.012
$data[35,2] = -124.52
$data[35,3] = 0
$data[35,4] = "top"
$data[36,0] = "C36"
$data[36,1] = 174.36500000000001
$data[36,2] = -116.48999999999999
$data[36,3] = 0
$data[36,4] = "top"
$data[37,0] = "C37"
$data[37,1] = 194.97999999999999
$data[37,2] = -117.88
$data[37,3] = 0
$data[37,4] = "top"
$data[38,0] = "C38"
$data[38,1] = 187.22999999999999
$data[38,2] = -114.33
$data[38,3] = 0
$data[38,4] = "top"
$data[39,0] = "C39"
$data[39,1] = 124.4965
$data[39,2] = -113.74369299999999
$data[39,3] = -90
$data[39,4] = "top"
$data[40,0] = "C40"
$data[40,1] = 108.9965
$data[40,2] = -113.943693
$data[40,3] = 0
$data[40,4] = "top"
$data[41,0] = "C41"
$data[41,1] = 116.8965
$data[41,2] = -105.343693
$data[41,3] = -90
$data[41,4] = "top"
$data[42,0] = "C42"
$data[42,1] = 73.496499999999997
$data[42,2] = -115.843693
$data[42,3] = 0
$data[42,4] = "top"
$data[43,0] = "C43"
$data[43,1] = 170.63
$data[43,2] = -98.340000000000003
$data[43,3] = 0
$data[43,4] = "top"
$data[44,0] = "C44"
$data[44,1] = 171
$data[44,2] = -88.890000000000001
$data[44,3] = 0
$data[44,4] = "top"
$data[45,0] = "C45"
$data[45,1] = 157.19999999999999
$data[45,2] = -98.340000000000003
$data[45,3] = 0
$data[45,4] = "top"
$data[46,0] = "C46"
$data[46,1] = 157.75999999999999
$data[46,2] = -88.890000000000001
$data[46,3] = 0
$data[46,4] = "top"
$data[47,0] = "D1"
$data[47,1] = 99.358500000000006
$data[47,2] = -120.681693
$data[47,3] = -90
$data[47,4] = "top"
$data[48,0] = "D2"
$data[48,1] = 88.436499999999995
$data[48,2] = -126.52369299999999
$data[48,3] = 90
$data[48,4] = "top"
$data[49,0] = "D3"
$data[49,1] = 85.896500000000003
$data[49,2] = -126.52369299999999
$data[49,3] = -90
$data[49,4] = "top"
$data[50,0] = "D4"
$data[50,1] = 99.485500000000002
$data[50,2] = -138.33469299999999
$data[50,3] = 90
$data[50,4] = "top"
$data[51,0] = "D5"
$data[51,1] = 88.436499999999995
$data[51,2] = -135.66769300000001
$data[51,3] = 90
$data[51,4] = "top"
$data[52,0] = "D6"
$data[52,1] = 85.896500000000003
$data[52,2] = -135.66769300000001
$data[52,3] = -90
$data[52,4] = "top"
$data[53,0] = "D7"
$data[53,1] = 63.596499999999999
$data[53,2] = -125.943693
$data[53,3] = 90
$data[53,4] = "top"
$data[54,0] = "D8"
$data[54,1] = 72.141499999999994
$data[54,2] = -143.94319300000001
$data[54,3] = -90
$data[54,4] = "top"
$data[55,0] = "D9"
$data[55,1] = 74.554500000000004
$data[55,2] = -143.94319300000001
$data[55,3] = -90
$data[55,4] = "top"
$data[56,0] = "D10"
$data[56,1] = 70.796499999999995
$data[56,2] = -124.543693
$data[56,3] = 90
$data[56,4] = "top"
$data[57,0] = "D11"
$data[57,1] = 170.78999999999999
$data[57,2] = -93.370000000000005
$data[57,3] = -90
$data[57,4] = "top"
$data[58,0] = "D12"
$data[58,1] = 157.28
$data[58,2] = -93.392499999999998
$data[58,3] = -90
$data[58,4] = "top"
$data[59,0] = "J1"
$data[59,1] = 64.236500000000007
$data[59,2] = -113.643693
$data[59,3] = -90
$data[59,4] = "top"
$data[60,0] = "J4"
$data[60,1] = 146.64993699999999
$data[60,2] = -68.622200000000007
$data[60,3] = 180
$data[60,4] = "top"
$data[61,0] = "J6"
$data[61,1] = 173.47
$data[61,2] = -129.02000000000001
$data[61,3] = 90
$data[61,4] = "top"
$data[62,0] = "J9"
$data[62,1] = 151.159559
$data[62,2] = -61.5
$data[62,3] = 180
$data[62,4] = "top"
$data[63,0] = "J10"
$data[63,1] = 190.08000000000001
$data[63,2] = -52.159999999999997
$data[63,3] = 90
$data[63,4] = "top"
$data[64,0] = "J11"
$data[64,1] = 67.519999999999996
$data[64,2] = -88.694999999999993
$data[64,3] = -90
$data[64,4] = "top"
$data[65,0] = "J12"
$data[65,1] = 60.223999999999997
$data[65,2] = -67.349999999999994
$data[65,3] = -90
$data[65,4] = "top"
$data[66,0] = "L1"
$data[66,1] = 92.754499999999993
$data[66,2] = -120.427693
$data[66,3] = 90
$data[66,4] = "top"
$data[67,0] = "L2"
$data[67,1] = 95.802499999999995
$data[67,2] = -139.09669299999999
$data[67,3] = 90
$data[67,4] = "top"
$data[68,0] = "L3"
$data[68,1] = 175.49000000000001
$data[68,2] = -98.340000000000003
$data[68,3] = 0
$data[68,4] = "top"
$data[69,0] = "L4"
$data[69,1] = 162
$data[69,2] = -98.340000000000003
$data[69,3] = 0
$data[69,4] = "top"
$data[70,0] = "Q1"
$data[70,1] = 67.046499999999995
$data[70,2] = -125.793693
$data[70,3] = 0
$data[70,4] = "top"
$data[71,0] = "Q3"
$data[71,1] = 86.896500000000003
$data[71,2] = -105.568693
$data[71,3] = 90
$data[71,4] = "top"
$data[72,0] = "Q4"
$data[72,1] = 86.796499999999995
$data[72,2] = -112.543693
$data[72,3] = -90
$data[72,4] = "top"
$data[73,0] = "R1"
$data[73,1] = 85.996499999999997
$data[73,2] = -109.068693
$data[73,3] = -90
$data[73,4] = "top"
$data[74,0] = "R2"
$data[74,1] = 72.6965
$data[74,2] = -119.643693
$data[74,3] = -90
$data[74,4] = "top"
$data[75,0] = "R3"
$data[75,1] = 72.246499999999997
$data[75,2] = -140.343693
$data[75,3] = -90
$data[75,4] = "top"
$data[76,0] = "R4"
$data[76,1] = 108.9965
$data[76,2] = -111.943693
$data[76,3] = 180
$data[76,4] = "top"
$data[77,0] = "R5"
$data[77,1] = 94.496499999999997
$data[77,2] = -106.74369299999999
$data[77,3] = 180
$data[77,4] = "top"
$data[78,0] = "R6"
$data[78,1] = 108.9965
$data[78,2] = -110.143693
$data[78,3] = 180
$data[78,4] = "top"
$data[79,0] = "R7"
$data[79,1] = 72.6965
$data[79,2] = -108.343693
$data[79,3] = 90
$data[79,4] = "top"
$data[80,0] = "R8"
$data[80,1] = 74.554500000000004
$data[80,2] = -140.343693
$data[80,3] = -90
$data[80,4] = "top"
$data[81,0] = "R9"
$data[81,1] = 71.4465
$data[81,2] = -128.49369300000001
$data[81,3] = 180
$data[81,4] = "top"
$data[82,0] = "R10"
$data[82,1] = 145.79650000000001
$data[82,2] = -138.943693
$data[82,3] = -90
$data[82,4] = "top"
$data[83,0] = "R11"
$data[83,1] = 87.596500000000006
$data[83,2] = -109.068693
$data[83,3] = 90
$data[83,4] = "top"
$data[84,0] = "R12"
$data[84,1] = 146.90000000000001
$data[84,2] = -75.650000000000006
$data[84,3] = 180
$data[84,4] = "top"
$data[85,0] = "R13"
$data[85,1] = 178.52000000000001
$data[85,2] = -91.599999999999994
$data[85,3] = 90
$data[85,4] = "top"
$data[86,0] = "R14"
$data[86,1] = 94.496499999999997
$data[86,2] = -109.943693
$data[86,3] = 0
$data[86,4] = "top"
$data[87,0] = "R15"
$data[87,1] = 108.9965
$data[87,2] = -108.443693
$data[87,3] = 0
$data[87,4] = "top"
$data[88,0] = "R16"
$data[88,1] = 94.496499999999997
$data[88,2] = -108.343693
$data[88,3] = 0
$data[88,4] = "top"
$data[89,0] = "R17"
$data[89,1] = 164.66999999999999
$data[89,2] = -91.599999999999994
$data[89,3] = 90
$data[89,4] = "top"
$data[90,0] = "R18"
$data[90,1] = 107.1777
$data[90,2] = -139.190293
$data[90,3] = 90
$data[90,4] = "top"
$data[91,0] = "R19"
$data[91,1] = 87.9285
$data[91,2] = -129.825693
$data[91,3] = 180
$data[91,4] = "top"
$data[92,0] = "R20"
$data[92,1] = 84.626499999999993
$data[92,2] = -129.825693
$data[92,3] = 180
$data[92,4] = "top"
$data[93,0] = "R21"
$data[93,1] = 67.717500000000001
$data[93,2] = -140.569693
$data[93,3] = -90
$data[93,4] = "top"
$data[94,0] = "R22"
$data[94,1] = 67.717500000000001
$data[94,2] = -143.843693
$data[94,3] = -90
$data[94,4] = "top"
$data[95,0] = "R23"
$data[95,1] = 87.9285
$data[95,2] = -132.111693
$data[95,3] = 180
$data[95,4] = "top"
$data[96,0] = "R24"
$data[96,1] = 84.626499999999993
$data[96,2] = -132.111693
$data[96,3] = 180
$data[96,4] = "top"
$data[97,0] = "R25"
$data[97,1] = 80.920000000000002
$data[97,2] = -90.719999999999999
$data[97,3] = 180
$data[97,4] = "top"
$data[98,0] = "R26"
$data[98,1] = 196.95500000000001
$data[98,2] = -129.27500000000001
$data[98,3] = 90
$data[98,4] = "top"
$data[99,0] = "R27"
$data[99,1] = 171.19
$data[99,2] = -122.59
$data[99,3] = 0
$data[99,4] = "top"
$data[100,0] = "R28"
$data[100,1] = 173.66499999999999
$data[100,2] = -118.69
$data[100,3] = 0
$data[100,4] = "top"
$data[101,0] = "R29"
$data[101,1] = 195.22999999999999
$data[101,2] = -115.848
$data[101,3] = 0
$data[101,4] = "top"
$data[102,0] = "R30"
$data[102,1] = 187.45500000000001
$data[102,2] = -112.43000000000001
$data[102,3] = 0
$data[102,4] = "top"
$data[103,0] = "R31"
$data[103,1] = 65.346500000000006
$data[103,2] = -122.24369299999999
$data[103,3] = 180
$data[103,4] = "top"
$data[104,0] = "R32"
$data[104,1] = 116.59650000000001
$data[104,2] = -111.818693
$data[104,3] = -90
$data[104,4] = "top"
$data[105,0] = "R33"
$data[105,1] = 108.9965
$data[105,2] = -106.693693
$data[105,3] = 180
$data[105,4] = "top"
$data[106,0] = "R34"
$data[106,1] = 80.909999999999997
$data[106,2] = -83.030000000000001
$data[106,3] = 0
$data[106,4] = "top"
$data[107,0] = "RN1"
$data[107,1] = 80.831500000000005
$data[107,2] = -86.570188999999999
$data[107,3] = 0
$data[107,4] = "top"
$data[108,0] = "TH1"
$data[108,1] = 109.2097
$data[108,2] = -139.190293
$data[108,3] = -90
$data[108,4] = "top"
$data[109,0] = "U1"
$data[109,1] = 113.8965
$data[109,2] = -126.74369299999999
$data[109,3] = -90
$data[109,4] = "top"
$data[110,0] = "U2"
$data[110,1] = 101.7465
$data[110,2] = -109.30619299999999
$data[110,3] = -90
$data[110,4] = "top"
$data[111,0] = "U3"
$data[111,1] = 100.12050000000001
$data[111,2] = -129.57169300000001
$data[111,3] = 180
$data[111,4] = "top"
$data[112,0] = "U4"
$data[112,1] = 73.496499999999997
$data[112,2] = -112.343693
$data[112,3] = 0
$data[112,4] = "top"
$data[113,0] = "U5"
$data[113,1] = 77.296499999999995
$data[113,2] = -122.143693
$data[113,3] = 180
$data[113,4] = "top"
$data[114,0] = "U6"
$data[114,1] = 174.66
$data[114,2] = -91.75
$data[114,3] = 90
$data[114,4] = "top"
$data[115,0] = "U7"
$data[115,1] = 120.2465
$data[115,2] = -111.943693
$data[115,3] = -90
$data[115,4] = "top"
$data[116,0] = "U8"
$data[116,1] = 131.82149999999999
$data[116,2] = -131.79369299999999
$data[116,3] = 180
$data[116,4] = "top"
$data[117,0] = "U9"
$data[117,1] = 79.996499999999997
$data[117,2] = -110.543693
$data[117,3] = 0
$data[117,4] = "top"
$data[118,0] = "U10"
$data[118,1] = 184.94999999999999
$data[118,2] = -120.70999999999999
$data[118,3] = 0
$data[118,4] = "top"
$data[119,0] = "U11"
$data[119,1] = 70.496499999999997
$data[119,2] = -133.48469299999999
$data[119,3] = -90
$data[119,4] = "top"
$data[120,0] = "U12"
$data[120,1] = 161.16999999999999
$data[120,2] = -91.75
$data[120,3] = 90
$data[120,4] = "top"
$data[121,0] = "Y1"
$data[121,1] = 120.1465
$data[121,2] = -105.893693
$data[121,3] = 0
$data[121,4] = "top"

# Write the entire table in one shot (covers existing rows 1-116 plus new rows 117-122)
$ws.Range("A1:E122").Value = $data

# Match the numeric display format already used for Mid X / Mid Y / Rotation columns.
# Rows 2-116 already carry this format from their existing style; only the newly
# appended rows (117-122) need it applied explicitly.
$ws.Range("B117:D122").NumberFormat = "0.000000"

# Give the newly-appended rows the same row height as the rest of the table
$ws.Range("A117:E122").RowHeight = 13.5

# Move the active selection to A1:E1, matching the saved workbook state
$ws.Range("A1:E1").Select()
